# Copy and paste the existing "Sector of Institution", "Reporting Location",
# and "Offense" fields (rows 2-10) twice -- into rows 11-19 and rows 20-28 --
# then replicate the "Date" operation (new sum2014/sum2015 labels) and the
# "Count" operation (new 2014/2015 counts) for each block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nudge the window position to match where the author left the Excel window
# (mirrors the xWindow/yWindow shift recorded in the workbook view).
$excel.ActiveWindow.Left = 7380
$excel.ActiveWindow.Top = 920

# ------------------------------------------------------------------
# Block 1: rows 2-10 copy/pasted into rows 11-19, dated "sum2014"
# ------------------------------------------------------------------
$ws.Range("A2:C10").Copy($ws.Range("A11"))

$ws.Range("D11:D19").Value = "sum2014"
$ws.Range("D2").Copy()
$ws.Range("D11:D19").PasteSpecial(-4122)

$counts2014 = @(622, 154, 10, 228, 7, 10, 24, 0, 3)
for ($i = 0; $i -lt $counts2014.Length; $i++) {
    $ws.Cells.Item(11 + $i, 5).Value = $counts2014[$i]
}
$ws.Range("E2").Copy()
$ws.Range("E11:E19").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Block 2: rows 2-10 copy/pasted into rows 20-28, dated "sum2015"
# ------------------------------------------------------------------
$ws.Range("A2:C10").Copy($ws.Range("A20"))

$ws.Range("D20:D28").Value = "sum2015"
$ws.Range("D2").Copy()
$ws.Range("D20:D28").PasteSpecial(-4122)

$counts2015 = @(723, 189, 15, 273, 15, 6, 23, 1, 4)
for ($i = 0; $i -lt $counts2015.Length; $i++) {
    $ws.Cells.Item(20 + $i, 5).Value = $counts2015[$i]
}
$ws.Range("E2").Copy()
$ws.Range("E20:E28").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Leave the selection where the author left off
$ws.Range("F21").Select()
